$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row before row 97, shifting rows 97:141 down to 98:142.
$ws.Rows.Item(97).Insert()

# Populate the newly inserted row 97 with the new weekly record.
$ws.Cells.Item(97, 1).Value = 9
$ws.Cells.Item(97, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(97, 3).Value = "Metropolitana"
$ws.Cells.Item(97, 4).Value = 44846
$ws.Cells.Item(97, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(97, 5).Value = 13
$ws.Cells.Item(97, 6).Value = 100112022
$ws.Cells.Item(97, 7).Value = "Arveja Verde"
$ws.Cells.Item(97, 8).Value = "Perfection"
$ws.Cells.Item(97, 9).Value = "Primera"
$ws.Cells.Item(97, 10).Value = 30
$ws.Cells.Item(97, 11).Value = 28000
$ws.Cells.Item(97, 12).Value = 28000
$ws.Cells.Item(97, 13).Value = 28000
$ws.Cells.Item(97, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(97, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(97, 16).Value = 1120
$ws.Cells.Item(97, 17).Value = 25
$ws.Cells.Item(97, 18).Value = "Hortaliza"
